# B3 - add/modify assessments
#
# 1) Delete the slide with sldId=262 (an empty, unused slide that was added
#    then removed again in the same editing session).
# 2) On slide 1 (sldId=256):
#    - shift the "Term 2" and "Term 3" rows further down to make room
#    - nudge the existing "Assessment" bar (Rectangle 6 / id 7) down slightly
#    - add a second "Assessment" bar (duplicate of the first) below it

$p = $ppt.ActivePresentation

# Shape.Left/.Top are backed by single-precision (32-bit) floats, and EMUs
# are recovered from them by truncation. Biasing the EMU->point conversion
# by half an EMU compensates for that truncation so the stored EMU value
# round-trips exactly back to the intended integer.
function EMU([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

function Move-Shape($shape, [double]$x, [double]$y) {
    $shape.Left = EMU $x
    $shape.Top = EMU $y
}

# ---------------------------------------------------------------------
# 1) Remove the extra slide (sldId 262) that was added then deleted again.
# ---------------------------------------------------------------------
$slideToDelete = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 262) {
        $slideToDelete = $candidate
    }
}
if ($slideToDelete -ne $null) {
    $slideToDelete.Delete()
}

# ---------------------------------------------------------------------
# 2) Update slide 1 (sldId 256).
# ---------------------------------------------------------------------
$s1 = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 256) {
        $s1 = $candidate
    }
}

# --- shift the Term 2 row down ---
Move-Shape ($s1.Shapes.Item("Rectangle 11")) 1793007 4198380
Move-Shape ($s1.Shapes.Item("TextBox 12"))   2196867 4303274
Move-Shape ($s1.Shapes.Item("TextBox 13"))   5644917 4272794
Move-Shape ($s1.Shapes.Item("TextBox 14"))   9260607 4271008

# --- shift the Term 3 row down ---
Move-Shape ($s1.Shapes.Item("Rectangle 15")) 1793007 5154928
Move-Shape ($s1.Shapes.Item("TextBox 16"))   2204487 5303756
Move-Shape ($s1.Shapes.Item("TextBox 17"))   5644917 5273276
Move-Shape ($s1.Shapes.Item("TextBox 18"))   9260607 5271490

# --- nudge the existing Assessment bar down slightly ---
$assessment1 = $s1.Shapes.Item("Rectangle 6")
Move-Shape $assessment1 3370347 3097172

# --- recreate PowerPoint's internal shape-name counter so the new shape we
#     keep ends up named "Rectangle 9" / id 10, matching what PowerPoint
#     produced (an intermediate shape "Rectangle 8" was created and removed
#     again during the original edit session) ---
$placeholder = $s1.Shapes.AddShape(1, 0, 0, 10, 10)
$placeholder.Delete()

# --- add the second Assessment bar as a duplicate of the first, so it
#     keeps identical formatting/style, then move & rename it ---
$dup = $assessment1.Duplicate()
$assessment2 = $dup.Item(1)
$assessment2.Name = "Rectangle 9"
Move-Shape $assessment2 3370347 3647237
